$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '71.347.90'
$ws.Range('E2').Value = '  +2.82%  '
$ws.Range('D3').Value = '3.588.43'
$ws.Range('E3').Value = '  +1.93%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '584.59'
$ws.Range('E5').Value = '  +2.49%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '186.79'
$ws.Range('E6').Value = '  +2.33%  '
$ws.Range('D7').Value = '3.577.48'
$ws.Range('E7').Value = '  +1.77%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.623'
$ws.Range('E8').Value = '  +1.37%  '
$ws.Range('E9').Value = '  +0.04%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.216'
$ws.Range('E10').Value = '  +15.86%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.654'
$ws.Range('E11').Value = '  +2.18%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '54.54'
$ws.Range('E12').Value = '  +1.85%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000319'
$ws.Range('E13').Value = '  +6.58%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.55'
$ws.Range('E14').Value = '  +1.12%  '
$ws.Range('D15').Value = '4.009.28'
$ws.Range('E15').Value = '  -1.99%  '
$ws.Range('D16').Value = '71.250.37'
$ws.Range('E16').Value = '  +2.86%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '19.28'
$ws.Range('E17').Value = '  +0.22%  '
$ws.Range('D18').Value = '3.574.10'
$ws.Range('E18').Value = '  +1.80%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.38'
$ws.Range('E19').Value = '  +0.53%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '565.76'
$ws.Range('E20').Value = '  +4.76%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.121'
$ws.Range('E21').Value = '  +0.64%  '
$ws.Range('E22').Value = '  -1.53%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '17.66'
$ws.Range('E23').Value = '  -9.70%  '
$ws.Range('E24').Value = '  +3.67%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.59'
$ws.Range('E25').Value = '  +5.59%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '94.91'
$ws.Range('E26').Value = '  +1.04%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.30'
$ws.Range('E27').Value = '  +1.76%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.95'
$ws.Range('E28').Value = '  +1.43%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.18'
$ws.Range('E29').Value = '  +1.12%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '32.64'
$ws.Range('E30').Value = '  +3.07%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.27'
$ws.Range('E31').Value = '  -0.92%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '12.32'
$ws.Range('E32').Value = '  -1.63%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.116'
$ws.Range('E33').Value = '  +1.34%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '64.20'
$ws.Range('E34').Value = '  -1.27%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.37'
$ws.Range('E35').Value = '  +8.26%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '547.91'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.419'
$ws.Range('E37').Value = '  +5.44%  '
$ws.Range('D38').Value = '0.0₃0810'
$ws.Range('E38').Value = '  +6.22%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '37.73'
$ws.Range('E39').Value = '  -0.86%  '
$ws.Range('E40').Value = '  +0.13%  '
$ws.Range('E41').Value = '  +7.90%  '
$ws.Range('D42').Value = '3.517.95'
$ws.Range('E42').Value = '  +11.85%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.48'
$ws.Range('E43').Value = '  +3.43%  '
$ws.Range('E44').Value = '  +1.83%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0448'
$ws.Range('E45').Value = '  +1.33%  '
$ws.Range('E46').Value = '  -0.48%  '
$ws.Range('E47').Value = '  -0.47%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.34'
$ws.Range('E48').Value = '  +1.62%  '
$ws.Range('E49').Value = '  +2.56%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.00'
$ws.Range('E50').Value = '  +0.17%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.46'
$ws.Range('E51').Value = '  +4.80%  '
